$p = $ppt.ActivePresentation

# 1) Update the cached "datetimeFigureOut" footer/date field text from
#    2024-02-04 to 2024-02-11 everywhere it appears: the slide master and
#    every slide layout's "Date Placeholder" shape.
$oldDate = "2024-02-04"
$newDate = "2024-02-11"

function Update-DateShapes($shapeColl) {
    for ($shpIdx = 1; $shpIdx -le $shapeColl.Count; $shpIdx++) {
        $sh = $shapeColl.Item($shpIdx)
        $tr = $sh.TextFrame.TextRange
        if ($tr.Text -eq $oldDate) {
            $tr.Text = $newDate
        }
    }
}

$master = $p.SlideMaster
Update-DateShapes $master.Shapes

for ($layoutIdx = 1; $layoutIdx -le $master.CustomLayouts.Count; $layoutIdx++) {
    $cl = $master.CustomLayouts.Item($layoutIdx)
    Update-DateShapes $cl.Shapes
}

# 2) Remove the first slide (the large "Landing" background/footer/chart
#    layout slide) leaving only the remaining (empty-background) slide.
$s = $p.Slides.Item(1)
$s.Delete()
